$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - new race data (US)
$ws.Cells.Item(7, 1).Value = 20200110
$ws.Cells.Item(7, 2).Value = "US"
$ws.Cells.Item(7, 3).Formula = "=(50+46+40+37+32+26+14+6)+(50+46+43+40+34+30+22+20)"
$ws.Cells.Item(7, 4).Formula = "=(37+34+22+14+40+0+50+43)+(24+46+32+40+22+0+0+6)"
$ws.Cells.Item(7, 5).Formula = "=(37+34+22+14+40+0+50+43)+(24+50+46+20+22+30+34+13)"
$ws.Cells.Item(7, 6).Formula = "=(46+37+14+40+50+26+43+0)+(50+46+40+20+22+30+34+43)"
$ws.Cells.Item(7, 7).Formula = "=(46+37+14+40+50+26+43+0)+(50+46+40+20+22+30+34+43)"
$ws.Cells.Item(7, 8).Value = 620

# Row 8 - new race data (TdS)
$ws.Cells.Item(8, 1).Value = 20200110
$ws.Cells.Item(8, 2).Value = "TdS"
$ws.Cells.Item(8, 3).Formula = "=(400+320+240+200+180+160+128+116)+(400+320+200+180+160+128+104+88)"
$ws.Cells.Item(8, 4).Formula = "=(400+200+180+128+160+240+116+320)+(144+400+320+180+72+128+88+80)"
$ws.Cells.Item(8, 5).Formula = "=(400+200+180+128+160+116+240+320)+(400+160+320+200+180+72+128+88)"
$ws.Cells.Item(8, 6).Formula = "=(400+200+180+128+240+116+160+320)+(400+160+320+200+180+72+128+88)"
$ws.Cells.Item(8, 7).Formula = "=(400+200+180+128+240+116+160+320)+(400+160+320+200+180+72+128+88)"
$ws.Cells.Item(8, 8).Value = 3380

# Update selection to match the final cursor position recorded in the diff
$ws.Range("G8").Select()
